# HealthEcho Table Description workbook edit
# - health_logs table: replace the "vital_signs / jsonb" row with four new
#   column rows (blood_glucose/String, blood_pressure/String,
#   oxygen_saturation/String, pulse_rate/int)
# - reminders table (and its rows) shift down by two rows to make room
# - view: scroll/zoom/selection updated to match the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the three extra data rows that follow row 27
# (shifts old row 30 "reminders" header + its rows down to 32..38,
#  and updates the dimension + mergeCells refs automatically)
$ws.Rows("30:31").Insert()

# Replace the old vital_signs/jsonb row with the new blood_glucose row,
# then fill in the three brand-new rows after it.
$ws.Range("A27").Value = "blood_glucose"
$ws.Range("B27").Value = "String"

$ws.Range("A28").Value = "blood_pressure"
$ws.Range("B28").Value = "String"

$ws.Range("A29").Value = "oxygen_saturation"
$ws.Range("B29").Value = "String"

$ws.Range("A30").Value = "pulse_rate"
$ws.Range("B30").Value = "int"

# Update window/view state: zoomed to 105%, scrolled so row 16 is the
# top-most visible row, with B30 as the active/selected cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
[void]$ws.Range("B30").Select()
$win.Zoom = 105
